$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update Kilinochchi center capacity value
$ws.Range("A13").Value = 156.12

# Update the active selection to match the saved workbook state
$ws.Range("A13").Select()
